# Update the cryptocurrency price/volume table on Sheet1 with the latest
# scraped values (price in column D, 1h volume % change in column E).
# Numeric-looking price strings are written with a leading apostrophe so
# Excel keeps them as text (preserving formats like trailing zeros,
# e.g. "1.00") instead of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.811.09'
$ws.Range("E2").Value = '  +3.08%  '

$ws.Range("D3").Value = '3.447.63'
$ws.Range("E3").Value = '  +2.00%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '''582.44'
$ws.Range("E5").Value = '  +2.30%  '

$ws.Range("D6").Value = '''146.87'
$ws.Range("E6").Value = '  +4.86%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '''0.478'
$ws.Range("E8").Value = '  +1.06%  '

$ws.Range("D9").Value = '''7.66'
$ws.Range("E9").Value = '  +0.63%  '

$ws.Range("E10").Value = '  +2.75%  '

$ws.Range("E11").Value = '  +1.70%  '

$ws.Range("D12").Value = '4.039.25'
$ws.Range("E12").Value = '  +2.07%  '

$ws.Range("D13").Value = '''29.16'
$ws.Range("E13").Value = '  +5.04%  '

$ws.Range("E14").Value = '  -0.65%  '

$ws.Range("D15").Value = '3.456.09'
$ws.Range("E15").Value = '  +2.15%  '

$ws.Range("E16").Value = '  +2.64%  '

$ws.Range("D17").Value = '62.821.86'

$ws.Range("D18").Value = '''6.24'
$ws.Range("E18").Value = '  +2.70%  '

$ws.Range("D19").Value = '''14.29'
$ws.Range("E19").Value = '  +5.61%  '

$ws.Range("E20").Value = '  +5.00%  '

$ws.Range("D21").Value = '''396.50'
$ws.Range("E21").Value = '  +4.16%  '

$ws.Range("D22").Value = '''0.565'
$ws.Range("E22").Value = '  +2.92%  '

$ws.Range("D23").Value = '''75.49'
$ws.Range("E23").Value = '  +0.11%  '

$ws.Range("E24").Value = '  +0.15%  '

$ws.Range("D25").Value = '''0.0000119'
$ws.Range("E25").Value = '  +4.33%  '

$ws.Range("D26").Value = '3.579.98'
$ws.Range("E26").Value = '  +1.75%  '

$ws.Range("E27").Value = '  -0.73%  '

$ws.Range("D28").Value = '''7.68'
$ws.Range("E28").Value = '  +6.64%  '

$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  +0.04%  '

$ws.Range("E30").Value = '  +3.00%  '

$ws.Range("D31").Value = '''2.15'
$ws.Range("E31").Value = '  +1.48%  '

$ws.Range("E32").Value = '  +5.17%  '

$ws.Range("E33").Value = '  +0.02%  '

$ws.Range("D34").Value = '''23.86'
$ws.Range("E34").Value = '  +2.90%  '

$ws.Range("D35").Value = '''5.33'
$ws.Range("E35").Value = '  +7.90%  '

$ws.Range("D36").Value = '''1.61'
$ws.Range("E36").Value = '  +11.20%  '

$ws.Range("D37").Value = '''7.07'
$ws.Range("E37").Value = '  +2.49%  '

$ws.Range("D38").Value = '''168.61'
$ws.Range("E38").Value = '  +1.31%  '

$ws.Range("D39").Value = '3.483.49'
$ws.Range("E39").Value = '  +2.00%  '

$ws.Range("D40").Value = '''30.07'
$ws.Range("E40").Value = '  +16.09%  '

$ws.Range("D41").Value = '''0.0770'
$ws.Range("E41").Value = '  +1.06%  '

$ws.Range("D42").Value = '''0.792'
$ws.Range("E42").Value = '  +1.64%  '

$ws.Range("D43").Value = '''4.48'
$ws.Range("E43").Value = '  +3.31%  '

$ws.Range("E44").Value = '  +5.38%  '

$ws.Range("E45").Value = '  +7.43%  '

$ws.Range("D46").Value = '2.528.50'
$ws.Range("E46").Value = '  +4.24%  '

$ws.Range("D47").Value = '''23.55'
$ws.Range("E47").Value = '  +3.68%  '

$ws.Range("E48").Value = '  +1.75%  '

$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").Value = '''2.20'
$ws.Range("E49").Value = '  +5.36%  '

$ws.Range("B50").Value = 'FirstDigitalUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D50").Value = '''1.00'
$ws.Range("E50").Value = '  +0.01%  '

$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = '''0.0268'
$ws.Range("E51").Value = '  +3.12%  '
